$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new blank row at row 17, pushing the existing rows 17-56 down to 18-57.
$ws.Rows("17:17").Insert()

# Populate the newly inserted row 17 with the new entry.
$ws.Range("R17").Value = "reconciliation tangedco"
$ws.Range("S17").Value = "2024-09-04 15:10:12"
